$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 978.36365
$ws.Range("I70").Value = 660.4
$ws.Range("J70").Value = 1243.3334
$ws.Range("K70").Value = 1981.2
$ws.Range("L70").Value = 3730.0002
$ws.Range("M70").Value = -1711.2
$ws.Range("N70").Value = -4270.0002
$ws.Range("H73").Value = 978.36365
$ws.Range("I73").Value = 660.4
$ws.Range("J73").Value = 1243.3334
$ws.Range("K73").Value = 1981.2
$ws.Range("L73").Value = 3730.0002
$ws.Range("M73").Value = -1045.2
$ws.Range("N73").Value = -5602.0002
$ws.Range("H116").Value = 2066.7778
$ws.Range("I116").Value = 2039
$ws.Range("J116").Value = 2101.5
$ws.Range("K116").Value = 2039
$ws.Range("L116").Value = 2101.5
$ws.Range("M116").Value = 1403
$ws.Range("N116").Value = -8985.5
$ws.Range("H132").Value = 4153.442
$ws.Range("I132").Value = 4332.2144
$ws.Range("J132").Value = 3819.7334
$ws.Range("K132").Value = 12996.6432
$ws.Range("L132").Value = 11459.2002
$ws.Range("M132").Value = -10466.6432
$ws.Range("N132").Value = -16519.2002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2428.6924
$ws.Range("I2").Value = 2448.8667
$ws.Range("J2").Value = 2401.182
$ws.Range("K2").Value = 2448.8667
$ws.Range("L2").Value = 2401.182
$ws.Range("M2").Value = -2335.8667
$ws.Range("N2").Value = -2627.182
$ws.Range("H32").Value = 24627.99
$ws.Range("I32").Value = 21288.486
$ws.Range("K32").Value = 21288.486
$ws.Range("M32").Value = -21001.486
$ws.Range("H45").Value = 998.2
$ws.Range("I45").Value = 872.75
$ws.Range("K45").Value = 872.75
$ws.Range("M45").Value = -495.75
$ws.Range("H61").Value = 83502020
$ws.Range("I61").Value = 125126330
$ws.Range("J61").Value = 253403.5
$ws.Range("K61").Value = 125126330
$ws.Range("L61").Value = 253403.5
$ws.Range("M61").Value = -125126118
$ws.Range("N61").Value = -253827.5
$ws.Range("H97").Value = 1736.3636
$ws.Range("I97").Value = 1711.1111
$ws.Range("J97").Value = 1850
$ws.Range("K97").Value = 1711.1111
$ws.Range("L97").Value = 1850
$ws.Range("M97").Value = -1215.1111
$ws.Range("N97").Value = -2842
$ws.Range("H116").Value = 2428.6924
$ws.Range("I116").Value = 2448.8667
$ws.Range("J116").Value = 2401.182
$ws.Range("K116").Value = 2448.8667
$ws.Range("L116").Value = 2401.182
$ws.Range("M116").Value = -154.8667
$ws.Range("N116").Value = -6989.182
$ws.Range("H132").Value = 40344.145
$ws.Range("I132").Value = 28468.5
$ws.Range("J132").Value = 66889.7
$ws.Range("K132").Value = 85405.5
$ws.Range("L132").Value = 200669.1
$ws.Range("M132").Value = -82875.5
$ws.Range("N132").Value = -205729.1
$ws.Range("H136").Value = 83502020
$ws.Range("I136").Value = 125126330
$ws.Range("J136").Value = 253403.5
$ws.Range("K136").Value = 375378990
$ws.Range("L136").Value = 760210.5
$ws.Range("M136").Value = -375376440
$ws.Range("N136").Value = -765310.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2428.6924
$ws.Range("I3").Value = 2448.8667
$ws.Range("J3").Value = 2401.182
$ws.Range("K3").Value = 2448.8667
$ws.Range("L3").Value = 2401.182
$ws.Range("M3").Value = -2334.8667
$ws.Range("N3").Value = -2629.182
$ws.Range("H99").Value = 882.7857
$ws.Range("I99").Value = 882.7857
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 882.7857
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 615.2143
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 4377.4707
$ws.Range("I107").Value = 3492.9092
$ws.Range("J107").Value = 5999.1665
$ws.Range("K107").Value = 3492.9092
$ws.Range("L107").Value = 5999.1665
$ws.Range("M107").Value = -1572.9092
$ws.Range("N107").Value = -9839.166499999999
$ws.Range("H132").Value = 49231.74
$ws.Range("J132").Value = 49231.74
$ws.Range("L132").Value = 49231.74
$ws.Range("N132").Value = -59351.74
$ws.Range("H134").Value = 2869.3142
$ws.Range("I134").Value = 3000
$ws.Range("J134").Value = 2428.25
$ws.Range("K134").Value = 9000
$ws.Range("L134").Value = 7284.75
$ws.Range("M134").Value = -6465
$ws.Range("N134").Value = -12354.75
$ws.Range("H135").Value = 45310.645
$ws.Range("J135").Value = 45310.645
$ws.Range("L135").Value = 45310.645
$ws.Range("N135").Value = -55450.645
$ws.Range("H141").Value = 54590
$ws.Range("J141").Value = 54590
$ws.Range("L141").Value = 54590
$ws.Range("N141").Value = -64950
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 214499.16
$ws.Range("I31").Value = 49850.477
$ws.Range("J31").Value = 310544.22
$ws.Range("K31").Value = 49850.477
$ws.Range("L31").Value = 310544.22
$ws.Range("M31").Value = -49555.477
$ws.Range("N31").Value = -311134.22
$ws.Range("H34").Value = 214499.16
$ws.Range("I34").Value = 49850.477
$ws.Range("J34").Value = 310544.22
$ws.Range("K34").Value = 49850.477
$ws.Range("L34").Value = 310544.22
$ws.Range("M34").Value = -49648.477
$ws.Range("N34").Value = -310948.22
$ws.Range("H122").Value = 2340.1
$ws.Range("I122").Value = 2048.5
$ws.Range("J122").Value = 2631.7
$ws.Range("K122").Value = 6145.5
$ws.Range("L122").Value = 7895.099999999999
$ws.Range("M122").Value = -3695.5
$ws.Range("N122").Value = -12795.1
$ws.Range("H132").Value = 28244.342
$ws.Range("I132").Value = 1677
$ws.Range("J132").Value = 145899.72
$ws.Range("K132").Value = 5031
$ws.Range("L132").Value = 437699.16
$ws.Range("M132").Value = -2501
$ws.Range("N132").Value = -442759.16
$ws.Range("H134").Value = 77993.234
$ws.Range("I134").Value = 989
$ws.Range("J134").Value = 201200
$ws.Range("K134").Value = 2967
$ws.Range("L134").Value = 603600
$ws.Range("M134").Value = -432
$ws.Range("N134").Value = -608670
$ws.Range("H135").Value = 48720
$ws.Range("J135").Value = 48720
$ws.Range("L135").Value = 48720
$ws.Range("N135").Value = -58860
$ws.Range("H139").Value = 53200
$ws.Range("J139").Value = 53200
$ws.Range("L139").Value = 53200
$ws.Range("N139").Value = -63480
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 2720
$ws.Range("J106").Value = 2720
$ws.Range("L106").Value = 8160
$ws.Range("N106").Value = -10052
$ws.Range("H131").Value = 969.34375
$ws.Range("J131").Value = 1076.6296
$ws.Range("L131").Value = 3229.8888
$ws.Range("N131").Value = -13309.8888
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4775.9
$ws.Range("J80").Value = 4832.375
$ws.Range("L80").Value = 4832.375
$ws.Range("N80").Value = -6828.375
$ws.Range("H83").Value = 4775.9
$ws.Range("J83").Value = 4832.375
$ws.Range("L83").Value = 24161.875
$ws.Range("N83").Value = -34145.875
$ws.Range("H102").Value = 1638.7894
$ws.Range("I102").Value = 1642.8
$ws.Range("J102").Value = 1634.3334
$ws.Range("K102").Value = 1642.8
$ws.Range("L102").Value = 1634.3334
$ws.Range("M102").Value = -20.79999999999995
$ws.Range("N102").Value = -4878.3334
$ws.Range("H126").Value = 1475
$ws.Range("I126").Value = 1475
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4425
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1955
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 81524.67999999999
$ws.Range("I132").Value = 46732.684
$ws.Range("J132").Value = 336666
$ws.Range("K132").Value = 140198.052
$ws.Range("L132").Value = 1009998
$ws.Range("M132").Value = -137668.052
$ws.Range("N132").Value = -1015058
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 4000
$ws.Range("J24").Value = 4000
$ws.Range("L24").Value = 4000
$ws.Range("N24").Value = -4686
$ws.Range("H82").Value = 1358.909
$ws.Range("I82").Value = 1083.5
$ws.Range("K82").Value = 1083.5
$ws.Range("M82").Value = -722.5
$ws.Range("H85").Value = 1358.909
$ws.Range("I85").Value = 1083.5
$ws.Range("K85").Value = 1083.5
$ws.Range("M85").Value = 164.5
$ws.Range("H93").Value = 1106.2174
$ws.Range("I93").Value = 1106.2174
$ws.Range("K93").Value = 1106.2174
$ws.Range("M93").Value = 141.7826
$ws.Range("H100").Value = 1876.8
$ws.Range("I100").Value = 1741.7142
$ws.Range("J100").Value = 1995
$ws.Range("K100").Value = 1741.7142
$ws.Range("L100").Value = 1995
$ws.Range("M100").Value = -1200.7142
$ws.Range("N100").Value = -3077
$ws.Range("H122").Value = 3512.825
$ws.Range("I122").Value = 3778.1052
$ws.Range("J122").Value = 3272.8096
$ws.Range("K122").Value = 11334.3156
$ws.Range("L122").Value = 9818.4288
$ws.Range("M122").Value = -8884.3156
$ws.Range("N122").Value = -14718.4288
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 463.72974
$ws.Range("I107").Value = 377.48
$ws.Range("J107").Value = 643.4167
$ws.Range("K107").Value = 1132.44
$ws.Range("L107").Value = 1930.2501
$ws.Range("M107").Value = 787.5599999999999
$ws.Range("N107").Value = -5770.2501
$ws.Range("H113").Value = 634.25
$ws.Range("I113").Value = 808.85
$ws.Range("J113").Value = 416
$ws.Range("K113").Value = 2426.55
$ws.Range("L113").Value = 1248
$ws.Range("M113").Value = -256.5500000000002
$ws.Range("N113").Value = -5588
$ws.Range("H122").Value = 964.25
$ws.Range("I122").Value = 957.53845
$ws.Range("K122").Value = 2872.61535
$ws.Range("M122").Value = -422.61535
$ws.Range("H132").Value = 75419.516
$ws.Range("I132").Value = 59725.234
$ws.Range("J132").Value = 102099.8
$ws.Range("K132").Value = 179175.702
$ws.Range("L132").Value = 306299.4
$ws.Range("M132").Value = -176645.702
$ws.Range("N132").Value = -311359.4
$ws.Range("H136").Value = 43730.125
$ws.Range("I136").Value = 25347.342
$ws.Range("J136").Value = 151400.72
$ws.Range("K136").Value = 76042.026
$ws.Range("L136").Value = 454202.16
$ws.Range("M136").Value = -73492.026
$ws.Range("N136").Value = -459302.16
